# correction in sa algorithm and 746 logs
# Updates the "Fitness" column (C) on the active sheet so that the best-so-far
# fitness values reflect the corrected simulated-annealing run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-19 (Generation 0-17) -> 7632
$ws.Range("C2:C19").Value = 7632

# Rows 20-25 (Generation 18-23) -> 7610
$ws.Range("C20:C25").Value = 7610

# Rows 26-36 (Generation 24-34) -> 7295
$ws.Range("C26:C36").Value = 7295

# Rows 37-71 (Generation 35-69) -> 7293
$ws.Range("C37:C71").Value = 7293
